$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (was the 2022-01-13 "Primera" entry, becomes the 2021-01-27 "Primera" entry)
$ws.Range("D2").Value = 44223
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 3500
$ws.Range("O2").Value = 4000
$ws.Range("P2").Value = 3750
$ws.Range("S2").Value = 1875

# Row 3 updates (was the 2022-01-13 "Segunda" entry, becomes the 2021-01-27 "Segunda" entry)
$ws.Range("D3").Value = 44223
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 3000
$ws.Range("O3").Value = 3000
$ws.Range("P3").Value = 3000
$ws.Range("S3").Value = 1500

# Row 4 updates (was the 2021-01-27 "Primera" entry, becomes the 2022-01-13 "Primera" entry)
$ws.Range("D4").Value = 44574
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 6000
$ws.Range("O4").Value = 7000
$ws.Range("P4").Value = 6500
$ws.Range("S4").Value = 3250

# Row 5 updates (was the 2021-01-27 "Segunda" entry, becomes the 2022-01-13 "Segunda" entry)
$ws.Range("D5").Value = 44574
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 5000
$ws.Range("O5").Value = 5000
$ws.Range("P5").Value = 5000
$ws.Range("S5").Value = 2500
